$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("C8").Value = 99
$ws.Range("C9").Value = 24
$ws.Range("G9").Value = "'6144.00"
$ws.Range("C10").Value = 33
$ws.Range("G10").Value = "'15576.00"
$ws.Range("A11").Value = ""
$ws.Range("C11").Value = 73
$ws.Range("D11").Value = "'2.0"
$ws.Range("E11").Value = "Rewiring of 3/5 pin 6 amp. Light plug point with 1.5 sq. mm nominal size  FR PVC insulated unsheathed flexible copper conductor 1.1 kV grade  and 1.5 sq. mm nominal size FR PVC insulated unsheathed flexible copper earth conductor 1.1 kV grade(IS:694)   in recessed ISI marked MMS ( IS:9537 P - III ) virgin material  PVC conduit & it's  ISI marked (IS:3419-1988) accessories, 1.2 mm thick  MS box with earth terminal of required size,  6 A  switch, 3/5 pin 6 A socket, 3.0 mm thick ISI marked (IS:2036-1995) phenolic laminated sheet, Al.alloy / Cadmium plated iron/ brass  screws, cup washers, making connections, testing etc. as required.  For specification of copper  Conductor,  Phenolic Laminated sheet's & Electrical/ Wiring accessories refer Chapter E - 04, E - 05 & E - 07 For additional technical parameters of product / work refer Annexure 'A' attached with this BSR"
$ws.Range("F11").Value = 0
$ws.Range("G11").Value = "'0.00"
$ws.Range("A12").Value = "P. point"
$ws.Range("C12").Value = 50
$ws.Range("D12").Value = "'6"
$ws.Range("E12").Value = "On board"
$ws.Range("F12").Value = 136
$ws.Range("G12").Value = "'6800.00"
$ws.Range("C13").Value = 74
$ws.Range("G13").Value = "'1702.00"
$ws.Range("C14").Value = 47
$ws.Range("G14").Value = "'2350.00"
$ws.Range("C15").Value = 81
$ws.Range("G15").Value = "'2673.00"
$ws.Range("C16").Value = 5
$ws.Range("G16").Value = "'390.00"
$ws.Range("C17").Value = 49
$ws.Range("G17").Value = "'1470.00"
$ws.Range("C18").Value = 70
$ws.Range("D18").Value = "'9.0"
$ws.Range("E18").Value = "Providing & Fixing of IS 11037:1984  marked  non modular socket size flush type 180 watt rotary minimum 5 step fan regulator with CM/L no. printed and made out from industrial grade Polycarbonate or fire resistant ABS material including making connection testing etc. as required.All as per pre approved by Engineer in charge.  For additional technical parameters of products/ work  , refer   Annexure ""A"" attached with this BSR ."
$ws.Range("F18").Value = 219
$ws.Range("G18").Value = "'15330.00"
$ws.Range("C19").Value = 13
$ws.Range("G19").Value = "'3939.00"
$ws.Range("C20").Value = 1
$ws.Range("C21").Value = 55
$ws.Range("G21").Value = "'3080.00"
$ws.Range("C22").Value = 49
$ws.Range("C23").Value = 88
$ws.Range("D23").Value = "'19"
$ws.Range("E23").Value = "2 x 2.5 sq. mm. + 1x1.5sqmm"
$ws.Range("F23").Value = 81
$ws.Range("G23").Value = "'7128.00"
$ws.Range("A24").Value = "Mtr."
$ws.Range("C24").Value = 10
$ws.Range("D24").Value = "'20"
$ws.Range("E24").Value = "2 x 4.0 sq. mm. + 1 x 2.5 sq. mm."
$ws.Range("F24").Value = 122
$ws.Range("G24").Value = "'1220.00"
$ws.Range("A25").Value = "Set"
$ws.Range("C25").Value = 98
$ws.Range("D25").Value = "'13.0"
$ws.Range("E25").Value = "Plate Earthing  as per IS:3043 with Hot dipped G.I. Earth plate of size 600mm x 600mm x 6.0mm by embodying  3 to 4 mtr. below the ground level with 20  mm dia. G.I. 'B' class watering Pipe ,including all accessories like nut, bolts, reducer, nipple, wire meshed funnel, and Heavy duty weather proof poly-propylene earth pit chamber with lockable Jam free lid suitable for safe working load 5000 Kg or more of size Top Dia. 225 to 260 mm, Bottom Dia 300 to 350 mm. and Height  250 to 300 mm. and embodying the pipe  complete with alternate layers salt and coke/ charcoal, testing of earth resistance for value of 5 ohms or less  as required & must record by engineer in charge during site visit and ensure to enter in measurment book.All as per pre approved by Engineer in charge.  For additional technical parameters of products/ work  , refer   Annexure ""A"" attached with this BSR .   "
$ws.Range("F25").Value = 5733
$ws.Range("G25").Value = "'561834.00"
$ws.Range("A26").Value = ""
$ws.Range("C26").Value = 36
$ws.Range("D26").Value = "'14.0"
$ws.Range("E26").Value = "Supply & Laying following size earth wire in horizontal or vertical run in ground/surface/recess including riveting, soldering, saddles,  making connection with GI/Cu purity purity >95%  thimble etc. as required.  All as per pre approved by Engineer in charge.  For additional technical parameters of products/ work  , refer   Annexure ""A"" attached with this BSR .   "
$ws.Range("F26").Value = 0
$ws.Range("G26").Value = "'0.00"
$ws.Range("A27").Value = "Mtr."
$ws.Range("C27").Value = 78
$ws.Range("D27").Value = "'23"
$ws.Range("E27").Value = "8 SWG G.I. ( Hot Dipped  ) Wire "
$ws.Range("F27").Value = 20
$ws.Range("G27").Value = "'1560.00"
$ws.Range("A28").Value = ""
$ws.Range("C28").Value = 36
$ws.Range("D28").Value = "'15.0"
$ws.Range("E28").Value = "Providing & Fixing of  BEE  Star rated copper wounded double ball bearing capacitor start, aluminium body & Metallic  blade ceiling  fan  Conforming to all the performance requirements laid down in IS 374:2019 including all amendments, as applicable ; & Carry BIS licensing (i.e. ISI marking) with down rod up to 80 cm with secondary support safety cable ( steel rope) , cotter pin with 3 x 1.5 sq.mm pvc insulated flexible copper conductor making connection testing etc. as required. All as per pre approved by Engineer in charge.  For additional technical parameters of products/ work  , refer   Annexure ""A"" attached with this BSR ."
$ws.Range("F28").Value = 0
$ws.Range("G28").Value = "'0.00"
$ws.Range("A29").Value = "Each"
$ws.Range("C29").Value = 42
$ws.Range("D29").Value = "'25"
$ws.Range("E29").Value = "1200 mm Sweep BEE 1 Star rated (service value >=4.0 to < 4.5 )"
$ws.Range("F29").Value = 1890
$ws.Range("G29").Value = "'79380.00"
$ws.Range("C30").Value = 9
$ws.Range("G30").Value = "'4428.00"
$ws.Range("C31").Value = 29
$ws.Range("A32").Value = "Each"
$ws.Range("C32").Value = 6
$ws.Range("D32").Value = "'30"
$ws.Range("E32").Value = " 6 A to 32 A rating"
$ws.Range("F32").Value = 187
$ws.Range("G32").Value = "'1122.00"
$ws.Range("A33").Value = ""
$ws.Range("C33").Value = 90
$ws.Range("D33").Value = "'31"
$ws.Range("E33").Value = "Double pole MCB(With B/C curve tripping Characteristics)"
$ws.Range("F33").Value = 0
$ws.Range("G33").Value = "'0.00"
$ws.Range("A34").Value = "Each"
$ws.Range("C34").Value = 3
$ws.Range("D34").Value = "'32"
$ws.Range("E34").Value = " 50/63 A rating"
$ws.Range("F34").Value = 900
$ws.Range("G34").Value = "'2700.00"
$ws.Range("A35").Value = ""
$ws.Range("C35").Value = 69
$ws.Range("D35").Value = "'18.0"
$ws.Range("E35").Value = "Providing & Fixing of Recessed/surface mounting heavy duty horizontal type Double Door ( Metal / Glazed )Distribution board with Metal end box made out from Galvanized steel / CRCA sheet not less then 1.2 mm thick  conforming to IS-8623-1 & 3 /  IEC 61439- 1 & 3, powder painted complete with reversible door (for double door DB only )100 amp.  insulated copper bus bar/shorting link , copper neutral link, copper earth link , color coded interconnecting wire set  of suitable rating and din bar,masking sheet,  making internal DB  terminations with copper lugs, Ferrules,  detachable gland plate, including making connections, testing etc. as required. OEM shall have submit  NABL / CPRI / ERDA accrediated   lab type test reports  & All as per pre approved by Engineer in charge. For additional technical parameters of product / work refer Annexure 'A' attached with this BSR"
$ws.Range("F35").Value = 0
$ws.Range("G35").Value = "'0.00"
$ws.Range("C36").Value = 31
$ws.Range("D36").Value = "'34"
$ws.Range("E36").Value = "Metal door (single phase) IK-09 and IP-43 with Metal end box"
$ws.Range("A37").Value = "Each"
$ws.Range("C37").Value = 31
$ws.Range("D37").Value = "'35"
$ws.Range("E37").Value = "8 Way (8+2)"
$ws.Range("F37").Value = 2184
$ws.Range("G37").Value = "'67704.00"
$ws.Range("C38").Value = 59
$ws.Range("C39").Value = 60
$ws.Range("C40").Value = 85
$ws.Range("G42").Value = "'786530.00"
$ws.Range("H42").Value = "'786530.00"
$ws.Range("G44").Value = "'786530.00"
$ws.Range("H44").Value = "'786530.00"
